$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was reported for this feria; it belongs at the
# top of the historical log (row 10, right after the most recent existing
# entries), so push the existing rows 10:50 down to 11:51 and populate the
# freshly inserted row with the new observation.
$ws.Rows("10").Insert()

$ws.Range("A10").Value = 4
$ws.Range("B10").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C10").Value = "Los Lagos"
$ws.Range("D10").Value = 44959
$ws.Range("E10").Value = 10
$ws.Range("F10").Value = "Fruta"
$ws.Range("G10").Value = 100101
$ws.Range("H10").Value = "Berries"
$ws.Range("I10").Value = 100101001
$ws.Range("J10").Value = "Arándano (blue)"
$ws.Range("K10").Value = "Sin especificar"
$ws.Range("L10").Value = "Primera"
$ws.Range("M10").Value = 200
$ws.Range("N10").Value = 2000
$ws.Range("O10").Value = 2200
$ws.Range("P10").Value = 2100
$ws.Range("Q10").Value = "$/bandeja 2 kilos"
$ws.Range("R10").Value = "Provincia de Curicó"
$ws.Range("S10").Value = 1050
$ws.Range("T10").Value = 2
